# Update column M (rows 2-25) with corrected simulation values
# per commit: "Finalización capitulo 3 tesis y correcion simulación de tamaños"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1.201684507136474
    3  = 0.6656833054784236
    4  = 0.6463047589365307
    5  = 0.6580474635624795
    6  = 0.757009341167266
    7  = 1.498482037695104
    8  = 0.6442398222717846
    9  = 0.6413163752053981
    10 = 0.6768021015152874
    11 = 0.7967757783005831
    12 = 0.7289362570375295
    13 = 0.8198293695888442
    14 = 0.6240174537503569
    15 = 0.6578207325390386
    16 = 0.6226131253722249
    17 = 0.5887061987231393
    18 = 0.7513672989687394
    19 = 1.537609721552292
    20 = 0.5762730678032315
    21 = 0.8781938061620314
    22 = 0.87939616941898
    23 = 0.6088878033734711
    24 = 0.8143911596851082
    25 = 0.5893681226693055
}

foreach ($row in $values.Keys) {
    $ws.Range("M$row").Value = $values[$row]
}
